# tp sl for ce being changed
# Append 19 new rows (114-132) to Sheet6, mirroring the structure of the
# existing trade rows (A:Trigger_Level_High_Low .. K:Activation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pick up the date/time number format already used in column E (e.g. E113)
# so new rows in that column render the same way as the existing data.
$dateFmt = $ws.Range("E113").NumberFormat

# NOTE: nested array literals like @(@(1,2),@(3,4)) get flattened by this
# interpreter, so rows are appended one at a time using the unary comma
# operator to keep each row as its own array element.
$data = @()
$data += ,@(37000,"MARKET",38800,"PE",45660,38810,38700,10,2,5,0)
$data += ,@(37000,"MARKET",38800,"PE",45660,38810,38700,10,2,5,0)
$data += ,@(37000,"MARKET",38800,"PE",45660,38810,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(37000,"MARKET",39000,"PE",45660,39100,38700,10,2,5,0)
$data += ,@(38900,"LIMIT",38920,"PE",45660,39040,38900,10,2,5,0)
$data += ,@(38900,"LIMIT",38920,"PE",45660,39040,38900,10,2,5,0)
$data += ,@(38900,"LIMIT",38920,"CE",45660,39040,38900,10,2,5,0)
$data += ,@(38900,"LIMIT",38920,"CE",45660,39040,38900,10,2,5,0)
$data += ,@(38900,"LIMIT",38915,"CE",45660,39020,38900,10,2,5,0)
$data += ,@(38800,"LIMIT",39000,"PE",45660,39020,38900,10,2,5,0)

$startRow = 114
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $row = $data[$idx]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).NumberFormat = $dateFmt

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
}
